# "added minigame support to planner"
#
# Two sheets change:
#   - "challenges": row 2 tweaked (K2 0->1, new M2), two new rows (3,4) added
#   - "tasks": 11 new header columns (O1:Y1) + 8 new data rows (5..12)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: challenges
# ---------------------------------------------------------------------------
$ch = $wb.Worksheets.Item("challenges")

# Row 2 edits: K2 0 -> 1, insert M2 = 1 (L2/N2 unchanged)
$ch.Range("K2").Value = 1
$ch.Range("M2").Value = 1

# Row 3 (new)
$ch.Range("A3").Value = 17
$ch.Range("B3").Value = 1
$ch.Range("C3").Value = "TASKS_COLLECTION"
$ch.Range("D3").Value = "G1"
$ch.Range("E3").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/3ad4d1db-b854-45cb-bcef-59dbaee47f6e.jpeg"
$ch.Range("F3").Value = "Generated by AI"
$ch.Range("G3").NumberFormat = "@"
$ch.Range("G3").Value = "122"
$ch.Range("H3").NumberFormat = "yyyy-mm-dd hh:mm"
$ch.Range("H3").Value = 45658.25
$ch.Range("I3").NumberFormat = "yyyy-mm-dd hh:mm"
$ch.Range("I3").Value = 45839.25
$ch.Range("K3").Value = 0
$ch.Range("L3").Value = 30
$ch.Range("M3").Value = 2
$ch.Range("N3").Value = 10080
$ch.Range("O3").Value = 0

# Row 4 (new)
$ch.Range("A4").Value = 17
$ch.Range("B4").Value = 2
$ch.Range("C4").Value = "TASKS_COLLECTION"
$ch.Range("D4").Value = "G2"
$ch.Range("E4").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/3ad4d1db-b854-45cb-bcef-59dbaee47f6e.jpeg"
$ch.Range("F4").Value = "Generated by AI"
$ch.Range("G4").NumberFormat = "@"
$ch.Range("G4").Value = "122"
$ch.Range("H4").NumberFormat = "yyyy-mm-dd hh:mm"
$ch.Range("H4").Value = 45658.25
$ch.Range("I4").NumberFormat = "yyyy-mm-dd hh:mm"
$ch.Range("I4").Value = 45839.25
$ch.Range("K4").Value = 0
$ch.Range("L4").Value = 40
$ch.Range("N4").Value = 10080
$ch.Range("O4").Value = 1

# ---------------------------------------------------------------------------
# Sheet: tasks
# ---------------------------------------------------------------------------
$tk = $wb.Worksheets.Item("tasks")

# New header cells O1:Y1 -- text labels that happen to look numeric, keep the
# same bold/bordered header style already used by A1:N1 (copy from N1).
$headerStyleSource = $tk.Range("N1")
$newHeaders = [ordered]@{
    "O1" = "0"
    "P1" = "1"
    "Q1" = "2"
    "R1" = "3"
    "S1" = "0.1"
    "T1" = "1.1"
    "U1" = "2.1"
    "V1" = "0.2"
    "W1" = "1.2"
    "X1" = "2.2"
    "Y1" = "3.1"
}
foreach ($addr in $newHeaders.Keys) {
    $cell = $tk.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $newHeaders[$addr]
    $cell.Style = $headerStyleSource.Style()
}

# Row 5 (new)
$tk.Range("A5").Value = 1
$tk.Range("B5").Value = "Buy half heart"
$tk.Range("D5").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$tk.Range("H5").Value = 7
$tk.Range("I5").Value = "ConfusingArrowsData"
$tk.Range("J5").Value = "ConfusingArrowsData"
$tk.Range("K5").Value = 0
$tk.Range("L5").Value = "[MINIGAME_BUY_HALF_HEART, STRICTLY_GREATER, 0],[MINIGAMESTATE_ID, EQUAL, 2]"
$tk.Range("M5").Value = -5
$tk.Range("N5").Value = "GameBus Studio"

# Row 6 (new)
$tk.Range("A6").Value = 1
$tk.Range("B6").Value = "Score 10 points"
$tk.Range("D6").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$tk.Range("H6").Value = 7
$tk.Range("I6").Value = "ConfusingArrowsData"
$tk.Range("J6").Value = "ConfusingArrowsData"
$tk.Range("K6").Value = 0
$tk.Range("L6").Value = "[MINIGAME_SCORE, STRICTLY_GREATER, 10],[MINIGAMESTATE_ID, EQUAL, 2]"
$tk.Range("M6").Value = 15
$tk.Range("N6").Value = "GameBus Studio"

# Row 7 (new)
$tk.Range("A7").Value = 1
$tk.Range("B7").Value = "Walk 500 meters"
$tk.Range("D7").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$tk.Range("H7").Value = 7
$tk.Range("I7").Value = "WALK"
$tk.Range("J7").Value = "WALK"
$tk.Range("K7").Value = 0
$tk.Range("L7").Value = "[DISTANCE, STRICTLY_GREATER, 499]"
$tk.Range("M7").Value = 10
$tk.Range("N7").Value = "GameBus Studio"

# Row 8 (new)
$tk.Range("A8").Value = 1
$tk.Range("B8").Value = "Call a friend/family member"
$tk.Range("D8").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$tk.Range("G8").Value = 1
$tk.Range("H8").Value = 7
$tk.Range("I8").Value = "WALK"
$tk.Range("J8").Value = "WALK"
$tk.Range("K8").Value = 0
$tk.Range("L8").Value = " [SECRET, EQUAL, fdjklagas37]"
$tk.Range("M8").Value = 10
$tk.Range("N8").Value = "GameBus Studio"

# Row 9 (new)
$tk.Range("A9").Value = 2
$tk.Range("B9").Value = "Buy half heart"
$tk.Range("D9").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$tk.Range("H9").Value = 7
$tk.Range("I9").Value = "ConfusingArrowsData"
$tk.Range("J9").Value = "ConfusingArrowsData"
$tk.Range("K9").Value = 0
$tk.Range("L9").Value = "[MINIGAME_BUY_HALF_HEART, STRICTLY_GREATER, 0],[MINIGAMESTATE_ID, EQUAL, 3]"
$tk.Range("M9").Value = -5
$tk.Range("N9").Value = "GameBus Studio"

# Row 10 (new)
$tk.Range("A10").Value = 2
$tk.Range("B10").Value = "Score 20 points"
$tk.Range("D10").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$tk.Range("H10").Value = 7
$tk.Range("I10").Value = "ConfusingArrowsData"
$tk.Range("J10").Value = "ConfusingArrowsData"
$tk.Range("K10").Value = 0
$tk.Range("L10").Value = "[MINIGAME_SCORE, STRICTLY_GREATER, 20],[MINIGAMESTATE_ID, EQUAL, 3]"
$tk.Range("M10").Value = 15
$tk.Range("N10").Value = "GameBus Studio"

# Row 11 (new)
$tk.Range("A11").Value = 2
$tk.Range("B11").Value = "Walk 1500 meters"
$tk.Range("D11").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$tk.Range("H11").Value = 7
$tk.Range("I11").Value = "WALK"
$tk.Range("J11").Value = "WALK"
$tk.Range("K11").Value = 0
$tk.Range("L11").Value = "[DISTANCE, STRICTLY_GREATER, 499]"
$tk.Range("M11").Value = 20
$tk.Range("N11").Value = "GameBus Studio"

# Row 12 (new)
$tk.Range("A12").Value = 2
$tk.Range("B12").Value = "Call a friend/family member"
$tk.Range("D12").Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$tk.Range("G12").Value = 1
$tk.Range("H12").Value = 7
$tk.Range("I12").Value = "WALK"
$tk.Range("J12").Value = "WALK"
$tk.Range("K12").Value = 0
$tk.Range("L12").Value = " [SECRET, EQUAL, fdjklagas37]"
$tk.Range("M12").Value = 10
$tk.Range("N12").Value = "GameBus Studio"

Write-Output "done"
